# week3 - seesion 1
# Rename the existing sheet "Hoja1" -> "C", add a new sheet "C++" right
# after it, populate "C++" with the session-1 data, size the columns to
# match, and leave "C++" as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.ActiveSheet
$ws1.Name = "C"

$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "C++"

# --- data for the new "C++" sheet -----------------------------------
$ws2.Range("B1").Value = "Serial"

$ws2.Range("A2").Value = "PI"
$ws2.Range("B2").Value = 1592.43

$ws2.Range("A3").Value = "SUM"

$ws2.Range("A4").Value = "2VEC"

$ws2.Range("A5").Value = "COUNTING"
$ws2.Range("B5").Value = 617.827

$ws2.Range("A6").Value = "JULIA"
$ws2.Range("B6").Value = 1443.72

# --- column widths ----------------------------------------------------
# Sheet "C": column A and all columns from C onward are narrow, column B
# is a bit wider (matches the original author's column sizing).
$ws1.Range("A1:AMK1").EntireColumn.ColumnWidth = 7.6667
$ws1.Range("B1").EntireColumn.ColumnWidth = 10.05

# Sheet "C++": every column uses the same (slightly wider) width.
$ws2.Range("A1:AMK1").EntireColumn.ColumnWidth = 10.6667

# --- selection / active sheet -----------------------------------------
# "C" keeps its selection at A1 (no longer the tab shown on open) while
# "C++" becomes the active tab with B7 selected.
[void]$ws1.Range("A1").Select()
[void]$ws2.Range("B7").Select()
